# "Generate Report for Handoff"
#
# The localization-status report moved from "In Translation" to
# "Ready for handoff": the status cells and the "latest xliff generate /
# handoff" timestamps are refreshed, and the (now wider) Status columns
# are resized to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -----------------
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$overview.Range("F2").Value = "Ready for handoff"   # de-de status
$zhcn.Range("C2").Value     = "Ready for handoff"   # Status column
$dede.Range("C2").Value     = "Ready for handoff"   # Status column

# --- Timestamps refreshed for the new handoff -------------------------------
$overview.Range("G2").Value = "2016-09-03 07:02:30"  # Latest HO Xliff Generate Date
$dede.Range("H2").Value     = "2016-09-03 07:02:30"  # Latest Handoff Datetime (de-de)
$zhcn.Range("H2").Value     = "2016-09-03 07:02:25"  # Latest Handoff Datetime (zh-cn)

# --- Widen the Status columns so the longer text isn't truncated -----------
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332  # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332  # column F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth     = 16.333333333333332  # column C (Status)
$dede.Columns.Item(3).ColumnWidth     = 16.333333333333332  # column C (Status)
